$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix a couple of typos in the row labels (column A)
$ws.Range("A9").Value = "Extensions of remarks  "
$ws.Range("A17").Value = "     Measures passed, House joint resolutions  "
